{"js": "// The only substantive textual change in the target revision is the\n// insertion of a single space character between the end of the sentence\n// \"...as my agent for Power of Attorney for Health Care.\" and the\n// immediately following Jinja tag \"{% endif %}\" inside the paragraph\n// that contains the merge-field template text. (Everything else in the\n// diff is re-save \"noise\": extra XML namespaces Word 2016+ adds, new\n// w:proofErr spell-check wrappers around camel/dotted template tokens,\n// latentStyles additions, and the removal of the transient _GoBack\n// bookmark \u2014 none of these affect the document's visible content.)\n//\n// Anchor on the unique substring that identifies the exact spot:\n// \"Health Care.{% endif %}{% for person\" only occurs once in the body.\nconst needle = \"Health Care.{% endif %}{% for person\";\n\nconst results = context.document.body.search(needle, { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target anchor text not found: \" + needle);\n}\n\n// Replace the matched range with the same text plus the inserted space,\n// preserving everything else exactly as-is.\nconst replacement = \"Health Care. {% endif %}{% for person\";\nresults.items[0].insertText(replacement, \"Replace\");\nawait context.sync();\n", "ps1": "# The only substantive textual change in the target revision is the\n# insertion of a single space character between the end of the sentence\n# \"...as my agent for Power of Attorney for Health Care.\" and the\n# immediately following Jinja tag \"{% endif %}\" inside the paragraph\n# that contains the merge-field template text. (Everything else in the\n# diff is re-save \"noise\" from a newer Word build: extra XML namespaces,\n# new w:proofErr spell-check wrappers around camel/dotted template\n# tokens, latentStyles additions, and the removal of the transient\n# _GoBack bookmark -- none of these affect the document's visible\n# content, so they are not reproduced here.)\n#\n# Anchor on the unique substring that identifies the exact spot:\n# \"Health Care.{% endif %}{% for person\" only occurs once in the document.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"Health Care.{% endif %}{% for person\"\n$find.Replacement.Text = \"Health Care. {% endif %}{% for person\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
